# Weekly roll of the "Pina" price series: insert a new week at the top
# (row 399) and push every existing week down by one row, dropping the
# oldest week off the bottom as a brand-new last row (496).
#
# Only the weekly-varying columns shift: D (fecha), L (calidad),
# M (volumen), N (precio minimo), O (precio maximo), P (precio promedio),
# Q (unidad de medida), R (origen), S (precio unidad), T (unidades por caja).
# The descriptive columns (A,B,C,E,F,G,H,I,J,K) are identical for every
# row in this block and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 399
$lastRow  = 495
$newLastRow = 496

$shiftCols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# 1) Snapshot every value we will need before writing anything, keyed by
#    [row][col] so writes never clobber reads. Use Value2 (not Value) so
#    date-formatted cells come back as raw serial numbers rather than
#    COM DateTime objects - writing a DateTime back into a brand-new cell
#    (row 496 below) would make Excel stamp it with a fresh "m/d/yyyy"
#    style instead of inheriting the column's existing date style.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = @{}
    foreach ($c in $shiftCols) {
        $snapshot[$r][$c] = $ws.Cells.Item($r, $c).Value2()
    }
}

# Also snapshot the full old last row (A..T) - it becomes the new row 496
# verbatim.
$fullLastRow = @{}
for ($c = 1; $c -le 20; $c++) {
    $fullLastRow[$c] = $ws.Cells.Item($lastRow, $c).Value2()
}

# 2) Push rows 399..495 down to 400..496 for the shifting columns only
#    (i.e. new row r = old row r-1), working from the bottom up so we
#    never overwrite a value we still need to read (we already snapshotted
#    everything, but keep the safe order regardless).
for ($r = $lastRow; $r -ge $firstRow + 1; $r--) {
    foreach ($c in $shiftCols) {
        $ws.Cells.Item($r, $c).Value = $snapshot[$r - 1][$c]
    }
}

# 3) The brand-new last row (496) gets the old last row's (495) full
#    contents, copied verbatim into every column A..T.
for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item($newLastRow, $c).Value = $fullLastRow[$c]
}
# Column D carries the sheet's date format; brand-new cells default to
# the generic style, so restore the same number format the rest of the
# D column uses (matches the existing date cells' style, s="2").
$ws.Cells.Item($newLastRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 4) Row 399 becomes the new (most recent) week's data.
$ws.Cells.Item($firstRow, 4).Value  = 45244
$ws.Cells.Item($firstRow, 12).Value = "Primera"
$ws.Cells.Item($firstRow, 13).Value = 150
$ws.Cells.Item($firstRow, 14).Value = 26000
$ws.Cells.Item($firstRow, 15).Value = 26000
$ws.Cells.Item($firstRow, 16).Value = 26000
$ws.Cells.Item($firstRow, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item($firstRow, 18).Value = "Ecuador"
$ws.Cells.Item($firstRow, 19).Value = 2167
$ws.Cells.Item($firstRow, 20).Value = 12
